# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" sheet (between "总计" and "2021-Q4") with fund
#   holdings data.
# - Add a "2022-Q3" summary row to the "总计" sheet (pushing the existing
#   "2021-Q4" row down).
#
# The existing "2021-Q4" worksheet is recreated (after capturing/knowing its
# contents) so that sheet IDs come out in the same order as a fresh insert
# (2022-Q3 = sheetId 2, 2021-Q4 = sheetId 3) rather than reusing the old
# "2021-Q4" sheet object (which would keep sheetId 2 for itself).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)    # "总计"
$wsOld   = $wb.Worksheets.Item(2)    # "2021-Q4" (existing)

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the existing "2021-Q4" row down to row 3 (copying
#    its formatted index cell A2 -> A3), then write the new "2022-Q3" row
#    into row 2.
# ---------------------------------------------------------------------

$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.16

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.4

# ---------------------------------------------------------------------
# 2) Drop the existing "2021-Q4" worksheet and recreate both per-quarter
#    sheets from scratch, in order, so the relationship/sheet ids line up
#    the way they would for a freshly inserted sheet.
# ---------------------------------------------------------------------

$wsOld.Delete() | Out-Null

$wsQ3 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ3.Name = "2022-Q3"
$wsQ3.Outline.SummaryRow = 1
$wsQ3.Outline.SummaryColumn = 1
$wsQ3.PageSetup.LeftMargin = 54
$wsQ3.PageSetup.RightMargin = 54
$wsQ3.PageSetup.TopMargin = 72
$wsQ3.PageSetup.BottomMargin = 72
$wsQ3.PageSetup.HeaderMargin = 36
$wsQ3.PageSetup.FooterMargin = 36

$wsQ4 = $wb.Worksheets.Add($null, $wsQ3)
$wsQ4.Name = "2021-Q4"
$wsQ4.Outline.SummaryRow = 1
$wsQ4.Outline.SummaryColumn = 1
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 3) Populate "2022-Q3" with the new fund holdings data.
# ---------------------------------------------------------------------

# Header row (styled like the other per-quarter sheets).
$wsTotal.Range("A2").Copy($wsQ3.Range("B1:H1"))
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Row 2
$wsTotal.Range("A2").Copy($wsQ3.Range("A2"))
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "'000739"
$wsQ3.Range("C2").Value = "平安新鑫先锋混合A"
$wsQ3.Range("D2").Value = "'5.38"
$wsQ3.Range("E2").Value = "'84.94"
$wsQ3.Range("F2").Value = "'3.88"
$wsQ3.Range("G2").Value = "'0.2087"
$wsQ3.Range("H2").Value = 5

# Row 3
$wsTotal.Range("A2").Copy($wsQ3.Range("A3"))
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "'001515"
$wsQ3.Range("C3").Value = "平安新鑫先锋混合C"
$wsQ3.Range("D3").Value = "'3.15"
$wsQ3.Range("E3").Value = "'84.94"
$wsQ3.Range("F3").Value = "'3.88"
$wsQ3.Range("G3").Value = "'0.1222"
$wsQ3.Range("H3").Value = 5

# Row 4
$wsTotal.Range("A2").Copy($wsQ3.Range("A4"))
$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = "'011807"
$wsQ3.Range("C4").Value = "平安研究精选混合A"
$wsQ3.Range("D4").Value = "'1.12"
$wsQ3.Range("E4").Value = "'89.39"
$wsQ3.Range("F4").Value = "'3.85"
$wsQ3.Range("G4").Value = "'0.0431"
$wsQ3.Range("H4").Value = 5

# Row 5
$wsTotal.Range("A2").Copy($wsQ3.Range("A5"))
$wsQ3.Range("A5").Value = 3
$wsQ3.Range("B5").Value = "'011808"
$wsQ3.Range("C5").Value = "平安研究精选混合C"
$wsQ3.Range("D5").Value = "'0.62"
$wsQ3.Range("E5").Value = "'89.39"
$wsQ3.Range("F5").Value = "'3.85"
$wsQ3.Range("G5").Value = "'0.0239"
$wsQ3.Range("H5").Value = 5

# ---------------------------------------------------------------------
# 4) Recreate "2021-Q4" with its original (unchanged) fund holdings data.
# ---------------------------------------------------------------------

# Header row
$wsTotal.Range("A2").Copy($wsQ4.Range("B1:H1"))
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Row 2
$wsTotal.Range("A2").Copy($wsQ4.Range("A2"))
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").Value = "'233009"
$wsQ4.Range("C2").Value = "大摩多因子精选策略混合"
$wsQ4.Range("D2").Value = "'8.14"
$wsQ4.Range("E2").Value = "'93.86"
$wsQ4.Range("F2").Value = "'1.24"
$wsQ4.Range("G2").Value = "'0.1009"
$wsQ4.Range("H2").Value = 7

# Row 3
$wsTotal.Range("A2").Copy($wsQ4.Range("A3"))
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").Value = "'009246"
$wsQ4.Range("C3").Value = "摩根士丹利华鑫ESG量化先行混合"
$wsQ4.Range("D3").Value = "'4.04"
$wsQ4.Range("E3").Value = "'93.45"
$wsQ4.Range("F3").Value = "'1.50"
$wsQ4.Range("G3").Value = "'0.0606"
$wsQ4.Range("H3").Value = 6
